# The deck's main theme (applied through the slide master -> ppt/theme/theme1.xml,
# which is also the presentation-level theme relationship) is switched from the
# "Integral" / "Red Violet" palette to the stock "Office Theme" palette - the
# palette that used to live only on the notes-master theme part. The font
# scheme and format scheme (fills/lines/effects) are already byte-identical
# between the two themes, so only the twelve theme colors need to change.

function HexToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# Index order matches the DrawingML <a:clrScheme> child order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRGB $officeTheme[$i - 1]
}

# Also rename the design/theme to match (no-op on hosts that treat this as
# read-only, harmless otherwise).
try {
    $p.Designs.Item(1).Name = "Office Theme"
} catch {
}

Write-Output "Theme palette swapped to Office Theme"
